$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension -> measure renames, and "aragon" column switches from
# its own iaest-dimension to the generic sdmx-dimension:refArea
$ws.Range("A2").Value = "iaest-measure:numero-de-viajes-diarios"
$ws.Range("C2").Value = "iaest-measure:situacion-preferente"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:lugar-trabajo-o-estudio"

# Row 3: these columns are now curated as "medida" instead of "dim"
$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: datatype changes from skos:Concept to xsd:int for measure columns,
# and the "aragon" column now points at a Comunidad URI instead of a concept
$ws.Range("A4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (mapping file references) is no longer needed and is removed
$ws.Rows(5).Delete()
